$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.050.17'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.16%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.607.35'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.65%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.44%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  +3.87%  '

$ws.Range("E9").Value = '  -0.53%  '

$ws.Range("E10").Value = '  -1.76%  '

$ws.Range("E11").Value = '  +5.29%  '

$ws.Range("E12").Value = '  -0.92%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.063.55'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.07%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '58.953.27'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.18%  '

$ws.Range("E15").Value = '  -2.42%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.605.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.54%  '

$ws.Range("E17").Value = '  -2.10%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.47'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '338.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.81%  '

$ws.Range("E20").Value = '  -2.24%  '

$ws.Range("E21").Value = '  -0.80%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.69'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.428'
$ws.Range("D24").Style = "Normal"

$ws.Range("E25").Value = '  -0.40%  '

$ws.Range("E26").Value = '  -2.16%  '

$ws.Range("E27").Value = '  -1.73%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0759'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.89%  '

$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("E30").Value = '  +1.38%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.84%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '154.14'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.57%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.96'
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.896'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.49%  '

$ws.Range("E36").Value = '  +5.01%  '

$ws.Range("E37").Value = '  -0.58%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '37.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.89%  '

$ws.Range("E39").Value = '  +0.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '283.94'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.85%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.997'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.18%  '

$ws.Range("E43").Value = '  -0.92%  '

$ws.Range("E44").Value = '  +0.57%  '

$ws.Range("E47").Value = '  +0.67%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.69'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.59%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.948.34'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '117.40'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.13'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.87%  '

